$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.649.85'
$ws.Range('E2').Value = '  -0.82%  '
$ws.Range('D3').Value = '2.434.42'
$ws.Range('E3').Value = '  -1.51%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '''505.34'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -2.80%  '
$ws.Range('D6').Value = '''129.05'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -2.55%  '
$ws.Range('E7').Value = '  -0.19%  '
$ws.Range('E8').Value = '  -1.40%  '
$ws.Range('D9').Value = '2.449.35'
$ws.Range('E9').Value = '  -0.97%  '
$ws.Range('E10').Value = '  -0.22%  '
$ws.Range('D11').Value = '''0.0954'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -3.92%  '
$ws.Range('E12').Value = '  -3.33%  '
$ws.Range('D13').Value = '''0.330'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -3.69%  '
$ws.Range('D14').Value = '2.869.02'
$ws.Range('E14').Value = '  -1.42%  '
$ws.Range('D15').Value = '57.584.77'
$ws.Range('E15').Value = '  -0.83%  '
$ws.Range('E16').Value = '  -1.32%  '
$ws.Range('D17').Value = '''0.0000133'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -2.82%  '
$ws.Range('D18').Value = '2.442.90'
$ws.Range('E18').Value = '  -1.25%  '
$ws.Range('D19').Value = '''10.46'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -3.63%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').Value = '''314.92'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -1.78%  '
$ws.Range('B21').Value = 'Polkadot'
$ws.Range('C21').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D21').Value = '''4.10'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -1.70%  '
$ws.Range('D22').Value = '''0.999'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('E23').Value = '  -1.54%  '
$ws.Range('D24').Value = '''63.40'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -1.62%  '
$ws.Range('D25').Value = '''0.407'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.25%  '
$ws.Range('D26').Value = '''1.00'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('E27').Value = '  -1.03%  '
$ws.Range('D28').Value = '''7.22'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -2.51%  '
$ws.Range('D29').Value = '''170.30'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +2.84%  '
$ws.Range('D30').Value = '0.0₃0723'
$ws.Range('E30').Value = '  -3.83%  '
$ws.Range('D31').Value = '''6.21'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -2.62%  '
$ws.Range('E32').Value = '  -2.81%  '
$ws.Range('E33').Value = '  +0.11%  '
$ws.Range('E35').Value = '  -0.15%  '
$ws.Range('D36').Value = '''17.74'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -2.20%  '
$ws.Range('E37').Value = '  -4.90%  '
$ws.Range('D38').Value = '''3.92'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -1.75%  '
$ws.Range('D39').Value = '''36.37'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -0.34%  '
$ws.Range('E40').Value = '  -2.82%  '
$ws.Range('D41').Value = '''0.756'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -4.92%  '
$ws.Range('D42').Value = '''270.59'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -2.22%  '
$ws.Range('E43').Value = '  -2.82%  '
$ws.Range('D44').Value = '''4.98'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -1.24%  '
$ws.Range('D45').Value = '''0.581'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -2.25%  '
$ws.Range('E46').Value = '  +0.23%  '
$ws.Range('D47').Value = '''119.55'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -5.05%  '
$ws.Range('D48').Value = '''0.0486'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -1.48%  '
$ws.Range('E49').Value = '  -4.00%  '
$ws.Range('E50').Value = '  -2.48%  '
$ws.Range('D51').Value = '''16.58'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -3.08%  '
